# Regenerate save_data to use K instead of Strike#, updating the
# previously-computed "K" column (G) values for rows 2-32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 4
    4  = 1
    5  = 2
    6  = 2
    7  = 5
    8  = 4
    9  = 5
    10 = 2
    11 = 1
    12 = 7
    13 = 3
    14 = 6
    15 = 5
    16 = 1
    17 = 4
    18 = 3
    19 = 2
    20 = 1
    21 = 6
    22 = 2
    23 = 4
    24 = 3
    25 = 1
    26 = 2
    27 = 5
    28 = 4
    29 = 3
    30 = 2
    31 = 2
    32 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
